$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2700
$ws.Range("J40").Value = 2900
$ws.Range("L40").Value = 2900
$ws.Range("N40").Value = -3250

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 38960028
$ws.Range("I116").Value = 22822190
$ws.Range("K116").Value = 22822190
$ws.Range("M116").Value = -22818748

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2603.3215
$ws.Range("I132").Value = 2178.196
$ws.Range("K132").Value = 6534.588
$ws.Range("M132").Value = -4004.588

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4005.3447
$ws.Range("J138").Value = 4500.625
$ws.Range("L138").Value = 13501.875
$ws.Range("N138").Value = -23781.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 754.4286
$ws.Range("J5").Value = 1399
$ws.Range("L5").Value = 1399
$ws.Range("N5").Value = -1623

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12555.964
$ws.Range("I32").Value = 8602.608
$ws.Range("K32").Value = 8602.608
$ws.Range("M32").Value = -8315.608

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 429203.38
$ws.Range("I45").Value = 619272.25
$ws.Range("J45").Value = 1548.5
$ws.Range("K45").Value = 619272.25
$ws.Range("L45").Value = 1548.5
$ws.Range("M45").Value = -618895.25
$ws.Range("N45").Value = -2302.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 559133.3
$ws.Range("I61").Value = 3242.8572
$ws.Range("J61").Value = 2504750
$ws.Range("K61").Value = 3242.8572
$ws.Range("L61").Value = 2504750
$ws.Range("M61").Value = -3030.8572
$ws.Range("N61").Value = -2505174

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5659
$ws.Range("I74").Value = 4886.357
$ws.Range("K74").Value = 4886.357
$ws.Range("M74").Value = -4012.357

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 5659
$ws.Range("I77").Value = 4886.357
$ws.Range("K77").Value = 24431.785
$ws.Range("M77").Value = -20063.785

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 4623370.5
$ws.Range("I102").Value = 5907101.5
$ws.Range("J102").Value = 1939
$ws.Range("K102").Value = 5907101.5
$ws.Range("L102").Value = 1939
$ws.Range("M102").Value = -5905479.5
$ws.Range("N102").Value = -5183

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 83336664
$ws.Range("I110").Value = 125003010
$ws.Range("K110").Value = 125003010
$ws.Range("M110").Value = -125000965

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3796.16
$ws.Range("I122").Value = 2422.2144
$ws.Range("J122").Value = 5544.8184
$ws.Range("K122").Value = 7266.6432
$ws.Range("L122").Value = 16634.4552
$ws.Range("M122").Value = -4816.6432
$ws.Range("N122").Value = -21534.4552

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 559133.3
$ws.Range("I136").Value = 3242.8572
$ws.Range("J136").Value = 2504750
$ws.Range("K136").Value = 9728.571599999999
$ws.Range("L136").Value = 7514250
$ws.Range("M136").Value = -7178.571599999999
$ws.Range("N136").Value = -7519350

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 754.4286
$ws.Range("J4").Value = 1399
$ws.Range("L4").Value = 1399
$ws.Range("N4").Value = -1629

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3880.0667
$ws.Range("I86").Value = 2022.4445
$ws.Range("K86").Value = 2022.4445
$ws.Range("M86").Value = -899.4445000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3880.0667
$ws.Range("I89").Value = 2022.4445
$ws.Range("K89").Value = 10112.2225
$ws.Range("M89").Value = -4496.2225

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 6397.375
$ws.Range("I99").Value = 5036.8
$ws.Range("K99").Value = 5036.8
$ws.Range("M99").Value = -3538.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2500
$ws.Range("I105").Value = 2500
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2500
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -753
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H128").Value = 4900
$ws.Range("I128").Value = 4900
$ws.Range("K128").Value = 14700
$ws.Range("M128").Value = -12210

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2784.361
$ws.Range("I134").Value = 2421.625
$ws.Range("K134").Value = 7264.875
$ws.Range("M134").Value = -4729.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3446.5938
$ws.Range("I31").Value = 2337.1052
$ws.Range("J31").Value = 5068.154
$ws.Range("K31").Value = 2337.1052
$ws.Range("L31").Value = 5068.154
$ws.Range("M31").Value = -2042.1052
$ws.Range("N31").Value = -5658.154

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3446.5938
$ws.Range("I34").Value = 2337.1052
$ws.Range("J34").Value = 5068.154
$ws.Range("K34").Value = 2337.1052
$ws.Range("L34").Value = 5068.154
$ws.Range("M34").Value = -2135.1052
$ws.Range("N34").Value = -5472.154

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 326241.62
$ws.Range("I58").Value = 2141.6667
$ws.Range("K58").Value = 2141.6667
$ws.Range("M58").Value = -1938.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1149
$ws.Range("I94").Value = 2610
$ws.Range("J94").Value = 905.5
$ws.Range("K94").Value = 2610
$ws.Range("L94").Value = 905.5
$ws.Range("M94").Value = -2159
$ws.Range("N94").Value = -1807.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 326241.62
$ws.Range("I136").Value = 2141.6667
$ws.Range("K136").Value = 6425.000100000001
$ws.Range("M136").Value = -3875.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 5586
$ws.Range("I114").Value = 3206.75
$ws.Range("J114").Value = 7965.25
$ws.Range("K114").Value = 9620.25
$ws.Range("L114").Value = 23895.75
$ws.Range("M114").Value = -6366.25
$ws.Range("N114").Value = -30403.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1028.3334
$ws.Range("I117").Value = 827.6
$ws.Range("K117").Value = 2482.8
$ws.Range("M117").Value = 959.1999999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1965.45
$ws.Range("I129").Value = 516.53845
$ws.Range("J129").Value = 4656.2856
$ws.Range("K129").Value = 1549.61535
$ws.Range("L129").Value = 13968.8568
$ws.Range("M129").Value = 3450.38465
$ws.Range("N129").Value = -23968.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 268164.6
$ws.Range("I132").Value = 404538.38
$ws.Range("J132").Value = 5907.385
$ws.Range("K132").Value = 1213615.14
$ws.Range("L132").Value = 17722.155
$ws.Range("M132").Value = -1211085.14
$ws.Range("N132").Value = -22782.155

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 708.75
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 708.75
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 5998.3335
$ws.Range("I32").Value = 6500
$ws.Range("J32").Value = 4995
$ws.Range("K32").Value = 6500
$ws.Range("L32").Value = 4995
$ws.Range("M32").Value = -6183
$ws.Range("N32").Value = -5629

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 129747.52
$ws.Range("I132").Value = 230656.69
$ws.Range("K132").Value = 691970.0700000001
$ws.Range("M132").Value = -689440.0700000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 7110.364
$ws.Range("I136").Value = 7180.579
$ws.Range("K136").Value = 21541.737
$ws.Range("M136").Value = -18991.737

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 24909.455
$ws.Range("I51").Value = 23500
$ws.Range("J51").Value = 25222.666
$ws.Range("K51").Value = 23500
$ws.Range("L51").Value = 25222.666
$ws.Range("M51").Value = -22990
$ws.Range("N51").Value = -26242.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1915.9166
$ws.Range("I126").Value = 1899.8334
$ws.Range("J126").Value = 1932
$ws.Range("K126").Value = 5699.5002
$ws.Range("L126").Value = 5796
$ws.Range("M126").Value = -3229.5002
$ws.Range("N126").Value = -10736

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 500733.38
$ws.Range("I132").Value = 806349.9
$ws.Range("K132").Value = 2419049.7
$ws.Range("M132").Value = -2416519.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 8096.091
$ws.Range("I136").Value = 8665.143
$ws.Range("K136").Value = 25995.429
$ws.Range("M136").Value = -23445.429
